# ---------------------------------------------------------------------------
# "added opf analysis, clean code, fixed bugs"
#
# The `dtypes` sheet documents the dtype of every column of every pandapower
# network table. This edit adds the metadata rows describing two new tables
# that were added to the workbook's analysis: `pwl_cost` and `poly_cost`
# (OPF cost tables). They are inserted right before the existing
# `bus_geodata` dtype rows, which therefore move down by 12 rows.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dtypes")

# 1. Push the existing bus_geodata dtype rows (218:220) down to (230:232),
#    carrying their formatting (the bordered/bold style on column A) along
#    with the values.
$ws.Range("A218:D220").Copy($ws.Range("A230:D232"))

# 2. Fill in the 12 new rows (218:229) describing the pwl_cost and
#    poly_cost tables.
$newRows = @(
    @("pwl_cost",  "power_type", "object"),
    @("pwl_cost",  "element",    "uint32"),
    @("pwl_cost",  "et",         "object"),
    @("pwl_cost",  "points",     "object"),
    @("poly_cost", "element",    "uint32"),
    @("poly_cost", "et",         "object"),
    @("poly_cost", "cp0_eur",          "float64"),
    @("poly_cost", "cp1_eur_per_mw",   "float64"),
    @("poly_cost", "cp2_eur_per_mw2",  "float64"),
    @("poly_cost", "cq0_eur",          "float64"),
    @("poly_cost", "cq1_eur_per_mvar", "float64"),
    @("poly_cost", "cq2_eur_per_mvar2","float64")
)

$startRow = 218
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value2 = $newRows[$i][0]
    $ws.Cells.Item($r, 3).Value2 = $newRows[$i][1]
    $ws.Cells.Item($r, 4).Value2 = $newRows[$i][2]
}

# 3. The new rows' column-A cells need the same bordered / bold / centered
#    formatting as every other row-index cell in column A. Clone it from an
#    existing styled cell (formats only) instead of rebuilding the format
#    property-by-property, which keeps the style table tidy (no duplicate
#    style entries).
$endRow = $startRow + $newRows.Count - 1   # 229
$ws.Range("A2").Copy()
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 4. Column A holds the running index (row number - 2) for every data row;
#    re-stamp it for every row from the first inserted row through the end
#    of the (now longer) table.
$lastRow = 232
for ($r = $startRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 2
}

# 5. Refresh the sheet's dimension / selection to reflect the new extent
#    (Excel records where the user was last looking when it saved).
$ws.Activate()
$ws.Range("R236").Select()

Write-Host "dtypes sheet updated: now $($ws.UsedRange.Rows.Count) rows"
